$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Control 12
$ws.Range("D2").Value = 0.9990943661158863
$ws.Range("E2").Value = 0.9990943661158863

# Row 3 - Control 18
$ws.Range("D3").Value = 0.01904853835483634
$ws.Range("E3").Value = 0.01904853835483634

# Row 4 - Control 34
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = [double]"8.642147004695193E-06"
$ws.Range("E4").Value = [double]"8.642147004695193E-06"

# Row 5 - Control 42
$ws.Range("D5").Value = [double]"3.333634417049986E-10"
$ws.Range("E5").Value = [double]"3.333634417049986E-10"

# Row 6 - Control 21
$ws.Range("D6").Value = [double]"5.383408987573162E-13"
$ws.Range("E6").Value = [double]"5.383408987573162E-13"

# Row 7 - MDD 35
$ws.Range("D7").Value = [double]"6.030909187427842E-11"
$ws.Range("E7").Value = 0.9999999999396909

# Row 8 - MDD 22
$ws.Range("D8").Value = 0.9999814204325361
$ws.Range("E8").Value = [double]"1.857956746387224E-05"

# Row 9 - MDD 50
$ws.Range("D9").Value = 0.9999884657185047
$ws.Range("E9").Value = [double]"1.153428149525126E-05"

# Row 10 - MDD 45
$ws.Range("D10").Value = [double]"8.783295142770593E-23"

# Row 11 - MDD 28
$ws.Range("D11").Value = [double]"1.38607593141852E-05"
$ws.Range("E11").Value = 0.9999861392406858
$ws.Range("F11").Value = 9.253073692321777
$ws.Range("G11").Value = 0.6
